$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("A8").Value2 = "NICOLAI_1984"
$ws.Range("C8").Value2 = 10243.0
$ws.Range("D8").Value2 = 60.534147724820706
$ws.Range("A9").Value2 = "RAYMER"
$ws.Range("C9").Value2 = 6463.0
$ws.Range("D9").Value2 = 1.291828248122254
$ws.Range("A10").Value2 = "SADRAEY"
$ws.Range("C10").Value2 = 6396.0
$ws.Range("D10").Value2 = 0.24176597168341904
$ws.Range("A11").Value2 = "JENKINSON"
$ws.Range("C11").Value2 = 21031.0
$ws.Range("D11").Value2 = 229.6098468027633
$ws.Range("A12").Value2 = "KROO"
$ws.Range("C12").Value2 = 7149.0
$ws.Range("D12").Value2 = 12.043212153152714
$ws.Range("A13").Value2 = "TORENBEEK_1976"
$ws.Range("C13").Value2 = 10802.0
$ws.Range("D13").Value2 = 69.29511507600442
$ws.Range("A14").Value2 = "TORENBEEK_2013"
$ws.Range("C14").Value2 = 7891.0
$ws.Range("D14").Value2 = 23.67226005043056
$ws.Range("A15").Value2 = "ROSKAM"
$ws.Range("C15").Value2 = 14388.0
$ws.Range("D15").Value2 = 125.4969557224173

$ws = $wb.Worksheets.Item("WING")
$ws.Range("A8").Value2 = "TORENBEEK_1982"
$ws.Range("C8").Value2 = 6631.0
$ws.Range("D8").Value2 = -1.9577167689958312
$ws.Range("A9").Value2 = "RAYMER"
$ws.Range("C9").Value2 = 8394.0
$ws.Range("D9").Value2 = 24.1090220843084
$ws.Range("A10").Value2 = "KROO"
$ws.Range("C10").Value2 = 7561.0
$ws.Range("D10").Value2 = 11.792746721402883
$ws.Range("A11").Value2 = "TORENBEEK_2013"
$ws.Range("C11").Value2 = 6138.0
$ws.Range("D11").Value2 = -9.246940963368484

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("A8").Value2 = "RAYMER"
$ws.Range("C8").Value2 = 507.0
$ws.Range("D8").Value2 = -30.904403094809947
$ws.Range("A9").Value2 = "SADRAEY"
$ws.Range("C9").Value2 = 1040.0
$ws.Range("D9").Value2 = 41.734557754236
$ws.Range("A10").Value2 = "JENKINSON"
$ws.Range("C10").Value2 = 700.0
$ws.Range("D10").Value2 = -4.601739973110383
$ws.Range("A11").Value2 = "TORENBEEK_1976"
$ws.Range("C11").Value2 = 52.0
$ws.Range("D11").Value2 = -92.91327211228821
$ws.Range("A12").Value2 = "KROO"
$ws.Range("C12").Value2 = 737.0
$ws.Range("D12").Value2 = 0.4407394854537826
$ws.Range("A13").Value2 = "HOWE"
$ws.Range("C13").Value2 = 1415.0
$ws.Range("D13").Value2 = 92.84076848292686
$ws.Range("A14").Value2 = "ROSKAM"
$ws.Range("C14").Value2 = 1523.0
$ws.Range("D14").Value2 = 107.55935717278983
$ws.Range("A15").Value2 = "NICOLAI_2013"
$ws.Range("C15").Value2 = 399.0
$ws.Range("D15").Value2 = -45.62299178467292

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("A8").Value2 = "RAYMER"
$ws.Range("C8").Value2 = 180.0
$ws.Range("D8").Value2 = -75.46901885022838
$ws.Range("A9").Value2 = "SADRAEY"
$ws.Range("C9").Value2 = 749.0
$ws.Range("D9").Value2 = 2.07613822877189
$ws.Range("A10").Value2 = "JENKINSON"
$ws.Range("C10").Value2 = 502.0
$ws.Range("D10").Value2 = -31.58581923785916
$ws.Range("A11").Value2 = "TORENBEEK_1976"
$ws.Range("C11").Value2 = 124.0
$ws.Range("D11").Value2 = -83.10087965237955
$ws.Range("A12").Value2 = "KROO"
$ws.Range("C12").Value2 = 488.0
$ws.Range("D12").Value2 = -33.49378443839695
$ws.Range("A13").Value2 = "HOWE"
$ws.Range("C13").Value2 = 1145.0
$ws.Range("D13").Value2 = 56.04429675826945
$ws.Range("A14").Value2 = "ROSKAM"
$ws.Range("C14").Value2 = 1523.0
$ws.Range("D14").Value2 = 107.55935717278983

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("A10").Value2 = "KUNDU"
$ws.Range("C10").Value2 = 694.0
$ws.Range("D10").Value2 = 14.492257635279127
$ws.Range("A11").Value2 = "JENKINSON"
$ws.Range("C11").Value2 = 704.0
$ws.Range("D11").Value2 = 16.142001981608797
$ws.Range("A12").Value2 = "ROSKAM"
$ws.Range("C12").Value2 = 687.0
$ws.Range("D12").Value2 = 13.337436592848357
$ws.Range("A17").Value2 = "KUNDU"
$ws.Range("C17").Value2 = 694.0
$ws.Range("D17").Value2 = 14.492257635279127
$ws.Range("A18").Value2 = "JENKINSON"
$ws.Range("C18").Value2 = 704.0
$ws.Range("D18").Value2 = 16.142001981608797
$ws.Range("A19").Value2 = "ROSKAM"
$ws.Range("C19").Value2 = 687.0
$ws.Range("D19").Value2 = 13.337436592848357

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("A11").Value2 = "KUNDU"
$ws.Range("C11").Value2 = 3265.0
$ws.Range("D11").Value2 = 23.303482559712165
$ws.Range("A13").Value2 = "TORENBEEK_2013"
$ws.Range("C13").Value2 = 3457.0
$ws.Range("D13").Value2 = 30.55440710839968
$ws.Range("A18").Value2 = "KUNDU"
$ws.Range("C18").Value2 = 3265.0
$ws.Range("D18").Value2 = 23.303482559712165
$ws.Range("A20").Value2 = "TORENBEEK_2013"
$ws.Range("C20").Value2 = 3457.0
$ws.Range("D20").Value2 = 30.55440710839968
